# The workbook originally has these sheets (in order):
#   Parameters, Cities, Cars, Orders, SimpleOrders, Routes
# Target layout:
#   Parameters, Cities, Vehicles, Orders, Routes
#
# i.e.
#   - "Cars"        sheet is renamed to "Vehicles"
#   - "Orders"      sheet (the old one, with Id/City/Demand/.../Deliver_to) is removed
#   - "SimpleOrders" sheet is renamed to "Orders" (taking the old Orders' place)
#   - "Routes" sheet is left as-is (becomes the active/selected tab automatically
#     once the earlier "Orders" sheet is removed, same as the source workbook)

$wb = $excel.ActiveWorkbook

# Rename "Cars" -> "Vehicles"
$wsCars = $wb.Worksheets.Item("Cars")
$wsCars.Name = "Vehicles"

# Remove the old, detailed "Orders" sheet entirely
$wsOldOrders = $wb.Worksheets.Item("Orders")
$wsOldOrders.Delete() | Out-Null

# Rename "SimpleOrders" -> "Orders"
$wsSimpleOrders = $wb.Worksheets.Item("SimpleOrders")
$wsSimpleOrders.Name = "Orders"
